{"js": "// Update the date stamp and all 25 \"two-digit \u00d7 two-digit\" practice\n// answers in the worksheet table to the next day's generated set.\nconst replacements = [\n  [\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"],\n  [\"71\u00d759=4189\", \"40\u00d796=3840\"],\n  [\"92\u00d757=5244\", \"86\u00d715=1290\"],\n  [\"98\u00d783=8134\", \"77\u00d711=847\"],\n  [\"44\u00d770=3080\", \"38\u00d725=950\"],\n  [\"74\u00d772=5328\", \"85\u00d719=1615\"],\n  [\"73\u00d789=6497\", \"44\u00d766=2904\"],\n  [\"70\u00d716=1120\", \"55\u00d789=4895\"],\n  [\"41\u00d730=1230\", \"99\u00d735=3465\"],\n  [\"27\u00d729=783\", \"35\u00d712=420\"],\n  [\"86\u00d726=2236\", \"61\u00d734=2074\"],\n  [\"82\u00d772=5904\", \"99\u00d773=7227\"],\n  [\"64\u00d765=4160\", \"43\u00d751=2193\"],\n  [\"85\u00d767=5695\", \"76\u00d729=2204\"],\n  [\"88\u00d789=7832\", \"91\u00d739=3549\"],\n  [\"49\u00d747=2303\", \"25\u00d776=1900\"],\n  [\"94\u00d761=5734\", \"92\u00d743=3956\"],\n  [\"86\u00d740=3440\", \"52\u00d718=936\"],\n  [\"84\u00d741=3444\", \"70\u00d713=910\"],\n  [\"74\u00d786=6364\", \"58\u00d752=3016\"],\n  [\"60\u00d772=4320\", \"73\u00d793=6789\"],\n  [\"40\u00d739=1560\", \"16\u00d735=560\"],\n  [\"52\u00d717=884\", \"73\u00d754=3942\"],\n  [\"77\u00d790=6930\", \"17\u00d725=425\"],\n  [\"33\u00d726=858\", \"60\u00d748=2880\"],\n  [\"96\u00d791=8736\", \"27\u00d757=1539\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and all 25 \"two-digit x two-digit\" practice\n# answers in the worksheet table to the next day's generated set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"),\n    @(\"71\u00d759=4189\", \"40\u00d796=3840\"),\n    @(\"92\u00d757=5244\", \"86\u00d715=1290\"),\n    @(\"98\u00d783=8134\", \"77\u00d711=847\"),\n    @(\"44\u00d770=3080\", \"38\u00d725=950\"),\n    @(\"74\u00d772=5328\", \"85\u00d719=1615\"),\n    @(\"73\u00d789=6497\", \"44\u00d766=2904\"),\n    @(\"70\u00d716=1120\", \"55\u00d789=4895\"),\n    @(\"41\u00d730=1230\", \"99\u00d735=3465\"),\n    @(\"27\u00d729=783\", \"35\u00d712=420\"),\n    @(\"86\u00d726=2236\", \"61\u00d734=2074\"),\n    @(\"82\u00d772=5904\", \"99\u00d773=7227\"),\n    @(\"64\u00d765=4160\", \"43\u00d751=2193\"),\n    @(\"85\u00d767=5695\", \"76\u00d729=2204\"),\n    @(\"88\u00d789=7832\", \"91\u00d739=3549\"),\n    @(\"49\u00d747=2303\", \"25\u00d776=1900\"),\n    @(\"94\u00d761=5734\", \"92\u00d743=3956\"),\n    @(\"86\u00d740=3440\", \"52\u00d718=936\"),\n    @(\"84\u00d741=3444\", \"70\u00d713=910\"),\n    @(\"74\u00d786=6364\", \"58\u00d752=3016\"),\n    @(\"60\u00d772=4320\", \"73\u00d793=6789\"),\n    @(\"40\u00d739=1560\", \"16\u00d735=560\"),\n    @(\"52\u00d717=884\", \"73\u00d754=3942\"),\n    @(\"77\u00d790=6930\", \"17\u00d725=425\"),\n    @(\"33\u00d726=858\", \"60\u00d748=2880\"),\n    @(\"96\u00d791=8736\", \"27\u00d757=1539\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
